$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '60.880.37'
Set-TextValue 'E2' '  -0.91%  '
Set-TextValue 'D3' '3.389.75'
Set-TextValue 'E3' '  -1.33%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  -0.02%  '
Set-TextValue 'D5' '571.35'
Set-TextValue 'E5' '  -0.71%  '
Set-TextValue 'D6' '141.80'
Set-TextValue 'E6' '  -2.44%  '
Set-TextValue 'D7' '3.390.34'
Set-TextValue 'E7' '  -1.35%  '
Set-TextValue 'E8' '  +0.04%  '
Set-TextValue 'E9' '  -0.23%  '
Set-TextValue 'E10' '  -1.83%  '
Set-TextValue 'E11' '  -1.86%  '
Set-TextValue 'E12' '  +2.31%  '
Set-TextValue 'D13' '3.968.02'
Set-TextValue 'E13' '  -1.34%  '
Set-TextValue 'E14' '  +2.02%  '
Set-TextValue 'D15' '28.23'
Set-TextValue 'E15' '  +0.73%  '
Set-TextValue 'E16' '  -1.23%  '
Set-TextValue 'D17' '3.394.32'
Set-TextValue 'E17' '  -1.08%  '
Set-TextValue 'D18' '60.957.13'
Set-TextValue 'E18' '  -0.96%  '
Set-TextValue 'D19' '6.17'
Set-TextValue 'E19' '  -1.80%  '
Set-TextValue 'D20' '13.86'
Set-TextValue 'E20' '  -2.65%  '
Set-TextValue 'D21' '8.95'
Set-TextValue 'E21' '  -4.86%  '
Set-TextValue 'D22' '384.18'
Set-TextValue 'E22' '  -2.75%  '
Set-TextValue 'E23' '  -1.43%  '
Set-TextValue 'D24' '74.47'
Set-TextValue 'E24' '  +0.93%  '
Set-TextValue 'E25' '  +0.35%  '
Set-TextValue 'E26' '  -4.78%  '
Set-TextValue 'D27' '3.527.60'
Set-TextValue 'E27' '  -1.30%  '
Set-TextValue 'E28' '  -1.38%  '
Set-TextValue 'E29' '  -0.19%  '
Set-TextValue 'D30' '7.38'
Set-TextValue 'E30' '  -2.82%  '
Set-TextValue 'E31' '  -3.24%  '
Set-TextValue 'B32' 'PancakeSwap'
Set-TextValue 'C32' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D32' '2.14'
Set-TextValue 'E32' '  -1.61%  '
Set-TextValue 'B33' 'Fetch.AI'
Set-TextValue 'C33' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D33' '1.42'
Set-TextValue 'E33' '  -2.54%  '
Set-TextValue 'E34' '  +0.01%  '
Set-TextValue 'E35' '  -1.66%  '
Set-TextValue 'D36' '6.98'
Set-TextValue 'E36' '  -0.38%  '
Set-TextValue 'D37' '167.44'
Set-TextValue 'E37' '  -0.04%  '
Set-TextValue 'D38' '3.419.96'
Set-TextValue 'E38' '  -1.28%  '
Set-TextValue 'E39' '  -2.55%  '
Set-TextValue 'D40' '1.48'
Set-TextValue 'E40' '  -4.29%  '
Set-TextValue 'D41' '0.0774'
Set-TextValue 'E41' '  -1.22%  '
Set-TextValue 'D42' '27.54'
Set-TextValue 'E42' '  +2.06%  '
Set-TextValue 'B43' 'FirstDigitalUSD'
Set-TextValue 'C43' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D43' '1.00'
Set-TextValue 'E43' '  -0.04%  '
Set-TextValue 'B44' 'Mantle'
Set-TextValue 'C44' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D44' '0.780'
Set-TextValue 'E44' '  -2.40%  '
Set-TextValue 'D45' '42.12'
Set-TextValue 'E45' '  -0.44%  '
Set-TextValue 'E46' '  -1.38%  '
Set-TextValue 'D47' '1.67'
Set-TextValue 'E47' '  -3.40%  '
Set-TextValue 'D49' '2.477.29'
Set-TextValue 'E49' '  -4.68%  '
Set-TextValue 'D50' '6.82'
Set-TextValue 'E50' '  -1.40%  '
Set-TextValue 'D51' '23.09'
Set-TextValue 'E51' '  -0.53%  '
